# Commit: "Wed, Jul 01, 2020  3:06:21 PM"
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the custom "Table_0" style {135606C0-A387-4BB4-B1C1-E018362C55DB}
#    to the built-in PowerPoint table style
#    {D96C3968-5336-442F-ABED-221C95AC55F2}.
# 2) The presentation's colour theme is switched back from the
#    "Integral / Red Violet" design to the default "Office Theme" colours
#    (a Design-gallery theme change).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------

$tableSlides = @(14, 15, 16)
$newStyleId  = "{D96C3968-5336-442F-ABED-221C95AC55F2}"

foreach ($slideIndex in $tableSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the design back to the plain Office Theme colours ------------

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as OLE RGB() integers (R + G*256 + B*65536).
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
